$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Barrierefreiheit auch für den Administrationsbereich" -> append "?"
#    (keeps the new "?" as its own run with the same Arial/24 run formatting
#     as the paragraph mark, by splitting the paragraph and merging it back)
# ---------------------------------------------------------------------------
$pBereich = $d.Paragraphs(6)
$pBereich.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(7)
$newPara.Range.Text = "?"
$pBereichAgain = $d.Paragraphs(6)
$joinPos = $pBereichAgain.Range.End
$d.Range($joinPos - 1, $joinPos).Delete()

# ---------------------------------------------------------------------------
# 2) "Installationsanleitung: Was soll rein? (Am Dienstag fragen)"
#    -> drop the trailing " (Am Dienstag fragen)" run entirely
# ---------------------------------------------------------------------------
$pInstall = $d.Paragraphs(7)
$findInstall = $pInstall.Range.Find
$findInstall.ClearFormatting()
$findInstall.Execute(" (Am Dienstag fragen)", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# 3) "HTTPS und SSL" -> append " (Christian)?"
# ---------------------------------------------------------------------------
$pHttps = $d.Paragraphs(8)
$pHttps.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs(9)
$newPara2.Range.Text = " (Christian)?"
$pHttpsAgain = $d.Paragraphs(8)
$joinPos2 = $pHttpsAgain.Range.End
$d.Range($joinPos2 - 1, $joinPos2).Delete()

# ---------------------------------------------------------------------------
# 4) Insert two new bullet paragraphs after "HTTPS und SSL (Christian)?":
#       "Verknüpfung Lehrer Fach?"
#       "Sollen alle Mock Ups rein oder nur die Hauptfunktionen?"
#    New paragraphs inherit the KeinLeerraum / numId=2 list formatting from
#    the paragraph they are split from, matching the original bullet runs.
# ---------------------------------------------------------------------------
$pHttps2 = $d.Paragraphs(8)
$pHttps2.Range.InsertParagraphAfter()
$pLink = $d.Paragraphs(9)
$pLink.Range.Text = "Verknüpfung Lehrer Fach?"

$pLink2 = $d.Paragraphs(9)
$pLink2.Range.InsertParagraphAfter()
$pMock = $d.Paragraphs(10)
$pMock.Range.Text = "Sollen alle Mock Ups rein oder nur die Hauptfunktionen?"

# ---------------------------------------------------------------------------
# 5) Remove the two trailing paragraphs:
#       "Qualitätsanforderungen müssen überprüfbar sein"
#       "Woran erkennt man ein erfülltes Kriterium"
#    (the empty "KeinLeerraum" paragraph right before them is kept)
# ---------------------------------------------------------------------------
$total = $d.Paragraphs.Count
$pQuality = $d.Paragraphs($total - 1)
$pCriterion = $d.Paragraphs($total)
$d.Range($pQuality.Range.Start, $pCriterion.Range.End).Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
